$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Delete()
$ws.Range("H12").Select()
